$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 67 (Harris Interactive, 10/20, n=1705) ---
$ws.Range("A67").Value = 20
$ws.Range("B67").Value = 2021
$ws.Range("C67").Value = 7
$ws.Range("D67").Value = 10
$ws.Range("E67").Value = 17
$ws.Range("F67").Value = "harris"
$ws.Range("G67").Value = "online"
$ws.Range("H67").Value = "included"
$ws.Range("I67").Value = 1705
$ws.Range("J67").Value = 1
$ws.Range("K67").Value = 1
$ws.Range("L67").Value = 10
$ws.Range("M67").Value = 2
$ws.Range("N67").Value = 3
$ws.Range("O67").Value = 8
$ws.Range("P67").Value = 4
$ws.Range("Q67").Value = 23
$ws.Range("T67").Value = 14
$ws.Range("U67").Value = "T_0.5"
$ws.Range("V67").Value = 1
$ws.Range("W67").Value = 16
$ws.Range("X67").Value = 17
$ws.Range("Y67").Value = "T_0.5"
$ws.Range("AA67").Value = "T_0.5"

# --- New row 68 (Harris Interactive, 10/20, n=1685) ---
$ws.Range("A68").Value = 20
$ws.Range("B68").Value = 2021
$ws.Range("C68").Value = 7
$ws.Range("D68").Value = 10
$ws.Range("E68").Value = 17
$ws.Range("F68").Value = "harris"
$ws.Range("G68").Value = "online"
$ws.Range("H68").Value = "included"
$ws.Range("I68").Value = 1685
$ws.Range("J68").Value = 1
$ws.Range("K68").Value = 1
$ws.Range("L68").Value = 10
$ws.Range("M68").Value = 2
$ws.Range("N68").Value = 3
$ws.Range("O68").Value = 9
$ws.Range("P68").Value = 4
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = 10
$ws.Range("U68").Value = "T_0.5"
$ws.Range("V68").Value = 2
$ws.Range("W68").Value = 16
$ws.Range("X68").Value = 17
$ws.Range("Y68").Value = "T_0.5"
$ws.Range("AA68").Value = "T_0.5"

# --- New row 69 (Harris Interactive, 10/20, n=1705) ---
$ws.Range("A69").Value = 20
$ws.Range("B69").Value = 2021
$ws.Range("C69").Value = 7
$ws.Range("D69").Value = 10
$ws.Range("E69").Value = 17
$ws.Range("F69").Value = "harris"
$ws.Range("G69").Value = "online"
$ws.Range("H69").Value = "included"
$ws.Range("I69").Value = 1705
$ws.Range("J69").Value = 1
$ws.Range("K69").Value = 1
$ws.Range("L69").Value = 10
$ws.Range("M69").Value = 2
$ws.Range("N69").Value = 3
$ws.Range("O69").Value = 9
$ws.Range("P69").Value = 5
$ws.Range("Q69").Value = 25
$ws.Range("S69").Value = 8
$ws.Range("U69").Value = "T_0.5"
$ws.Range("V69").Value = 2
$ws.Range("W69").Value = 16
$ws.Range("X69").Value = 18
$ws.Range("Y69").Value = "T_0.5"
$ws.Range("AA69").Value = "T_0.5"

# --- Update view/selection to reflect the new bottom row ---
$excel.Goto($ws.Range("G1"), $false)
$ws.Range("Z69").Select()
